$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B233").Value = 63255
$ws.Range("F233").Value = 78
$ws.Range("G233").Value = 6396
$ws.Range("B234").Value = 57004
$ws.Range("F234").Value = 5
$ws.Range("G234").Value = 410
$ws.Range("F235").Value = 10
$ws.Range("G235").Value = 616.5
$ws.Range("B247").Value = 89604.12
$ws.Range("F263").Value = 100
$ws.Range("G263").Value = 6480
$ws.Range("F265").Value = 0
$ws.Range("G265").Value = 0
$ws.Range("F266").Value = 0
$ws.Range("G266").Value = 0
$ws.Range("F267").Value = 0
$ws.Range("G267").Value = 0
$ws.Range("F268").Value = 0
$ws.Range("G268").Value = 0
$ws.Range("F269").Value = 0
$ws.Range("G269").Value = 0
$ws.Range("B270").Value = 6526.49
$ws.Range("F273").Value = 19
$ws.Range("G273").Value = 1938.19
$ws.Range("B280").Value = 102883.72
$ws.Range("B322").Value = 48719
$ws.Range("C322").Value = "HIM-BABY CARE GIFT PACK (WW)1"
$ws.Range("D322").Value = 295.75
$ws.Range("E322").Value = 353.35
$ws.Range("F322").Value = -82
$ws.Range("G322").Value = -24251.5
$ws.Range("B323").Value = 66188
$ws.Range("C323").Value = "HIM-Baby Care Gift Pack (Ww)1"
$ws.Range("D323").Value = 315.8
$ws.Range("E323").Value = 377.31
$ws.Range("F323").Value = 35
$ws.Range("G323").Value = 11053
$ws.Range("F335").Value = 12
$ws.Range("G335").Value = 799.08
$ws.Range("F336").Value = 159
$ws.Range("G336").Value = 6754.32
$ws.Range("F355").Value = 32
$ws.Range("G355").Value = 3107.2
$ws.Range("B367").Value = 66194
$ws.Range("C367").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F367").Value = 35
$ws.Range("G367").Value = 2998.8
$ws.Range("B368").Value = 64983
$ws.Range("C368").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F368").Value = 6
$ws.Range("G368").Value = 514.08
$ws.Range("B369").Value = 64985
$ws.Range("C369").Value = "HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S"
$ws.Range("F369").Value = 13
$ws.Range("G369").Value = 1140.1
$ws.Range("B370").Value = 66196
$ws.Range("C370").Value = "HIM-Total Care Baby Pants Drapers-Xl-9S"
$ws.Range("F370").Value = 28
$ws.Range("G370").Value = 2455.6
$ws.Range("B372").Value = 142423.13
$ws.Range("B375").Value = 63565
$ws.Range("E375").Value = 109.19
$ws.Range("F375").Value = 60
$ws.Range("G375").Value = 6162.6
$ws.Range("B376").Value = 61610
$ws.Range("E376").Value = 122.71
$ws.Range("F376").Value = -58
$ws.Range("G376").Value = -5957.18
$ws.Range("B397").Value = 63560
$ws.Range("E397").Value = 134.87
$ws.Range("F397").Value = 1
$ws.Range("G397").Value = 126.86
$ws.Range("B398").Value = 60325
$ws.Range("E398").Value = 151.57
$ws.Range("F398").Value = -102
$ws.Range("G398").Value = -12939.72
$ws.Range("F457").Value = 39
$ws.Range("G457").Value = 29020.29
$ws.Range("B458").Value = 106236.89
$ws.Range("B548").Value = 65068
$ws.Range("E548").Value = 13.97
$ws.Range("F548").Value = 0
$ws.Range("G548").Value = 0
$ws.Range("B549").Value = 53602
$ws.Range("E549").Value = 15.69
$ws.Range("F549").Value = -232
$ws.Range("G549").Value = -3050.8
$ws.Range("B556").Value = 64922
$ws.Range("E556").Value = 20.98
$ws.Range("F556").Value = 0
$ws.Range("G556").Value = 0
$ws.Range("B557").Value = 45706
$ws.Range("E557").Value = 23.58
$ws.Range("F557").Value = -207
$ws.Range("G557").Value = -4084.11
$ws.Range("B564").Value = 45709
$ws.Range("E564").Value = 15.69
$ws.Range("F564").Value = -302
$ws.Range("G564").Value = -3971.3
$ws.Range("B565").Value = 64925
$ws.Range("E565").Value = 13.97
$ws.Range("F565").Value = 0
$ws.Range("G565").Value = 0
$ws.Range("B566").Value = 64919
$ws.Range("E566").Value = 27.97
$ws.Range("F566").Value = 0
$ws.Range("G566").Value = 0
$ws.Range("B567").Value = 45702
$ws.Range("E567").Value = 31.43
$ws.Range("F567").Value = -224
$ws.Range("G567").Value = -5891.2
$ws.Range("B569").Value = 53595
$ws.Range("E569").Value = 17.61
$ws.Range("F569").Value = -338
$ws.Range("G569").Value = -4978.74
$ws.Range("B570").Value = 65067
$ws.Range("E570").Value = 15.65
$ws.Range("F570").Value = 0
$ws.Range("G570").Value = 0
$ws.Range("F606").Value = 58
$ws.Range("G606").Value = 8428.559999999999
$ws.Range("F609").Value = 71
$ws.Range("G609").Value = 15778.33
$ws.Range("B612").Value = 128259.43
$ws.Range("B640").Value = 53319
$ws.Range("E640").Value = 310.64
$ws.Range("F640").Value = -6
$ws.Range("G640").Value = -1643.52
$ws.Range("B641").Value = 64810
$ws.Range("E641").Value = 291.22
$ws.Range("F641").Value = 2
$ws.Range("G641").Value = 547.84
$ws.Range("B659").Value = 60025
$ws.Range("E659").Value = 37.22
$ws.Range("F659").Value = -98
$ws.Range("G659").Value = -3217.34
$ws.Range("B660").Value = 64833
$ws.Range("E660").Value = 34.9
$ws.Range("F660").Value = 88
$ws.Range("G660").Value = 2889.04
$ws.Range("B669").Value = 60022
$ws.Range("E669").Value = 37.22
$ws.Range("F669").Value = -113
$ws.Range("G669").Value = -3709.79
$ws.Range("B670").Value = 64830
$ws.Range("E670").Value = 34.9
$ws.Range("F670").Value = 89
$ws.Range("G670").Value = 2921.87
$ws.Range("F702").Value = 49
$ws.Range("G702").Value = 2115.82
$ws.Range("B705").Value = 35964.32
$ws.Range("F827").Value = 37
$ws.Range("G827").Value = 13978.97
$ws.Range("B839").Value = 278378.29
$ws.Range("F890").Value = 1493
$ws.Range("G890").Value = 243523.23
$ws.Range("F892").Value = 52
$ws.Range("G892").Value = 14709.24
$ws.Range("B896").Value = 270029.04
$ws.Range("B941").Value = 3954180.08
$ws.Range("B942").Value = 3954180.08
